$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date value that was shifted forward by 10 days
# (from serial 45233 / 2023-11-03 to serial 45243 / 2023-11-13) for rows 2-32.
$newSerial = 45243

for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $newSerial
}
